# Generate Report for Handoff
# Rename the handed-off source file (old GUID -> new GUID) across all three
# sheets (Overview, zh-cn, de-de), refresh the "Latest HO Xliff Generate
# Date" / handoff timestamps, and rename the per-language xliff filenames.
#
# The hyperlinks on A2/B2 keep pointing at the same external URL (the
# commit/file reference on GitHub did not change), only their *displayed*
# text changes to reflect the new file name - so we delete + recreate the
# hyperlink on each sheet with the same Address but a new TextToDisplay.

$wb = $excel.ActiveWorkbook

$oldId = "46e19c79-0cb8-458e-83cf-68d8ba758c50"
$newId = "957017b5-0150-4d00-b914-918dcd26b2af"

$oldZhXlf = "$oldId.37636619180e1a395848851f9f141b10518df373.zh-cn.xlf"
$newZhXlf = "$newId.c3e58c6d068126f57cc1777d5468290609a10407.zh-cn.xlf"

$oldDeXlf = "$oldId.37636619180e1a395848851f9f141b10518df373.de-de.xlf"
$newDeXlf = "$newId.c3e58c6d068126f57cc1777d5468290609a10407.de-de.xlf"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a3d1d280dc1ed5ccd6734b052a43a0273dca3cf/e2e/$oldId.md"

# ---------------------------------------------------------------------
# Overview sheet: File Name (A2), Path And Name (B2, hyperlinked),
# Latest HO Xliff Generate Date (G2)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "e2e\$newId.md")

$wsOverview.Range("G2").Value = "2016-08-19 13:00:22"

# ---------------------------------------------------------------------
# zh-cn sheet: Source File Name (A2, hyperlinked), Latest Handoff File
# (G2), Latest Handoff Datetime (H2)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newId.md")

$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-19 13:00:03"

# ---------------------------------------------------------------------
# de-de sheet: Source File Name (A2, hyperlinked), Latest Handoff File
# (G2), Latest Handoff Datetime (H2)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newId.md")

$wsDeDe.Range("G2").Value = $newDeXlf
